# "Getting Form type improved"
# The sheet used to hold a literal form id ("form_1") next to the "TYPE"
# label. Replace it with a real, human readable form type, widen the
# column so the longer label fits (like Excel's AutoFit would), and leave
# the selection sitting on the cell that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "PERMANENT FOOD FACILITY"

# Widen column B so the new, longer value fits (matches AutoFit sizing).
$ws.Columns.Item(2).ColumnWidth = 25.7

$ws.Range("B1").Select()
